# The tutor originally supplied decimal "weight" values for the second
# exercise (Harok2). They should have been whole numbers (integers).
# This script rounds those inputs to the nearest integer in both places
# where they appear on the sheet (column B rows 4-23, and the duplicated
# column C rows 43-62 used by the regression table). Every other figure
# on the sheet (averages, sums, correlation, regression lines, the
# confidence-interval table B69:G79, and the charts that plot this data)
# is formula-driven and will recompute automatically once these raw
# inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Rounded integer weights (previously entered with long decimal tails).
$weights = @(67, 52, 78, 57, 109, 97, 58, 95, 98, 82, 68, 90, 85, 54, 104, 110, 75, 119, 94, 65)

# Column B, rows 4..23 (first listing of the weights)
for ($i = 0; $i -lt $weights.Length; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 2).Value = $weights[$i]
}

# Column C, rows 43..62 (second listing of the same weights, used by the
# regression / correlation computations further down the sheet)
for ($i = 0; $i -lt $weights.Length; $i++) {
    $row = 43 + $i
    $ws.Cells.Item($row, 3).Value = $weights[$i]
}

# Force a full recalculation so every dependent formula (averages, sums,
# SUMSQ/SUMPRODUCT, CORREL, the regression coefficients, the confidence
# bands in B69:G79, and the cached chart values) is refreshed to match
# the corrected inputs.
$excel.CalculateFullRebuild()

# Restore the saved scroll/selection position of the sheet, as recorded
# after the fix was made. (Best effort: this engine's Window object only
# exposes ScrollRow/ScrollColumn for the top-left visible cell.)
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
$ws.Range("J73").Select() | Out-Null
